$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.695.65'
$ws.Range("E2").Value = '  +1.14%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.704.01'
$ws.Range("E3").Value = '  -2.81%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '599.98'
$ws.Range("E5").Value = '  +0.77%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '167.66'
$ws.Range("E6").Value = '  -3.85%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.702.75'
$ws.Range("E7").Value = '  -2.81%  '

$ws.Range("E8").Value = '  -0.03%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.534'
$ws.Range("E9").Value = '  +1.30%  '

$ws.Range("E10").Value = '  +3.21%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.28'
$ws.Range("E11").Value = '  +0.04%  '

$ws.Range("E12").Value = '  -0.95%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '38.09'
$ws.Range("E13").Value = '  +0.20%  '

$ws.Range("E14").Value = '  -0.45%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.323.07'
$ws.Range("E15").Value = '  -2.79%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.702.82'
$ws.Range("E16").Value = '  -3.00%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '68.665.07'
$ws.Range("E17").Value = '  +0.89%  '

$ws.Range("E18").Value = '  +1.45%  '

$ws.Range("E19").Value = '  -0.49%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.24'
$ws.Range("E20").Value = '  +6.25%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '494.29'
$ws.Range("E21").Value = '  +0.98%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.20'
$ws.Range("E22").Value = '  +0.17%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.721'
$ws.Range("E23").Value = '  -1.17%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.47'
$ws.Range("E24").Value = '  -0.13%  '

$ws.Range("E25").Value = '  -3.75%  '

$ws.Range("E26").Value = '  +1.99%  '

$ws.Range("E27").Value = '  -0.78%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.06'
$ws.Range("E28").Value = '  -1.54%  '

$ws.Range("E29").Value = '  +0.11%  '

$ws.Range("E30").Value = '  +0.11%  '

$ws.Range("E31").Value = '  +1.98%  '

$ws.Range("E32").Value = '  -1.87%  '

$ws.Range("E33").Value = '  -3.72%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.841.38'
$ws.Range("E34").Value = '  -2.88%  '

$ws.Range("E35").Value = '  -0.65%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.644.24'
$ws.Range("E36").Value = '  -2.84%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  +0.02%  '

$ws.Range("E38").Value = '  -0.50%  '

$ws.Range("E39").Value = '  -0.36%  '

$ws.Range("E40").Value = '  -2.96%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.321'
$ws.Range("E41").Value = '  -1.04%  '

$ws.Range("E42").Value = '  +0.39%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '431.57'
$ws.Range("E43").Value = '  -3.52%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.97'
$ws.Range("E44").Value = '  -0.71%  '

$ws.Range("E45").Value = '  -2.38%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.37'
$ws.Range("E46").Value = '  +1.33%  '

$ws.Range("E47").Value = '  -0.02%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '40.23'
$ws.Range("E48").Value = '  -2.96%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '142.07'
$ws.Range("E49").Value = '  +2.53%  '

$ws.Range("B50").Value = 'VeChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0349'
$ws.Range("E50").Value = '  -0.36%  '

$ws.Range("B51").Value = 'Maker'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.735.75'
$ws.Range("E51").Value = '  -3.41%  '
